$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("A10").Value = "TestCase_F9"
$ws.Range("B10").Value = "OPQA-216"
$ws.Range("C10").Value = "Verify that user receives a notification when someone he is following user comments on a post"
$ws.Range("D10").Value = "Y"
$ws.Range("E10").Value = "PASS"

$ws.Range("A10:E10").Style = $ws.Range("A9:E9").Style

$ws.Range("D6").Select()
